# Export with no is_pref and no lev distance
#
# The sheet lists speaker variants per row. Previously, multiple spelling
# variants could be grouped under one shared "id" (column B), chosen via a
# Levenshtein-distance match against a "preferred" variant flagged with an
# "x" in column D (is_prefered). The new export drops that fuzzy grouping
# and the is_prefered flag entirely: every row now gets its own id derived
# directly from its own speaker_variant text, column D is always blank, and
# rows are re-ordered by (the row's own) id/name grouping.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New speaker_variant (column C) order for rows 2-15, after dropping the
# Levenshtein-based re-use of ids and re-grouping each variant with its own
# derived id.
$names = @(
    "AEmilius",
    "Fabius",
    "Marcelles",
    "Keiser",
    "Balbinus",
    "Calisthenes",
    "Faustina",
    "Valerius",
    "Vituria",
    "Julia",
    "Marcellus",
    "Baleinus",
    "Keizer",
    "Attilius"
)

$row = 2
foreach ($name in $names) {
    # id (column B) is now always derived from this row's own name, never
    # borrowed from a different (Levenshtein-nearest) preferred row.
    $id = "#" + $name.ToLower()

    $ws.Cells.Item($row, 2).Value = $id
    $ws.Cells.Item($row, 3).Value = $name
    # is_prefered (column D) is no longer emitted.
    $ws.Cells.Item($row, 4).Value = ""

    $row++
}
